$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 157.5051493644714

$ws.Range("A4").Value = 18813.4895
$ws.Range("B4").Value = 17525
$ws.Range("F4").Value = 6412.021
$ws.Range("G4").Value = 6455
